{"js": "// Update the worksheet date and the 25 two-digit-by-two-digit\n// multiplication prompts to the new values from the next day's sheet.\nconst replacements = [\n  [\"2025-12-21 Sunday\", \"2025-12-22 Monday\"],\n  [\"86\u00d777=\", \"38\u00d738=\"],\n  [\"54\u00d780=\", \"36\u00d770=\"],\n  [\"74\u00d793=\", \"26\u00d718=\"],\n  [\"43\u00d772=\", \"99\u00d725=\"],\n  [\"36\u00d757=\", \"89\u00d769=\"],\n  [\"78\u00d713=\", \"13\u00d731=\"],\n  [\"51\u00d753=\", \"68\u00d716=\"],\n  [\"21\u00d796=\", \"23\u00d745=\"],\n  [\"52\u00d788=\", \"43\u00d717=\"],\n  [\"99\u00d764=\", \"12\u00d731=\"],\n  [\"94\u00d747=\", \"61\u00d785=\"],\n  [\"29\u00d782=\", \"70\u00d761=\"],\n  [\"81\u00d743=\", \"39\u00d790=\"],\n  [\"72\u00d798=\", \"48\u00d745=\"],\n  [\"56\u00d760=\", \"55\u00d773=\"],\n  [\"18\u00d768=\", \"53\u00d737=\"],\n  [\"71\u00d751=\", \"73\u00d795=\"],\n  [\"52\u00d738=\", \"21\u00d780=\"],\n  [\"81\u00d780=\", \"35\u00d725=\"],\n  [\"84\u00d779=\", \"31\u00d757=\"],\n  [\"61\u00d743=\", \"54\u00d783=\"],\n  [\"34\u00d778=\", \"73\u00d735=\"],\n  [\"81\u00d737=\", \"20\u00d731=\"],\n  [\"84\u00d736=\", \"78\u00d723=\"],\n  [\"85\u00d764=\", \"20\u00d712=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first (and expected only) occurrence of each unique\n  // old value so numbers are never double-substituted.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 two-digit-by-two-digit\n# multiplication prompts to the new values from the next day's sheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-12-21 Sunday\"; New = \"2025-12-22 Monday\" },\n    @{ Old = \"86\u00d777=\"; New = \"38\u00d738=\" },\n    @{ Old = \"54\u00d780=\"; New = \"36\u00d770=\" },\n    @{ Old = \"74\u00d793=\"; New = \"26\u00d718=\" },\n    @{ Old = \"43\u00d772=\"; New = \"99\u00d725=\" },\n    @{ Old = \"36\u00d757=\"; New = \"89\u00d769=\" },\n    @{ Old = \"78\u00d713=\"; New = \"13\u00d731=\" },\n    @{ Old = \"51\u00d753=\"; New = \"68\u00d716=\" },\n    @{ Old = \"21\u00d796=\"; New = \"23\u00d745=\" },\n    @{ Old = \"52\u00d788=\"; New = \"43\u00d717=\" },\n    @{ Old = \"99\u00d764=\"; New = \"12\u00d731=\" },\n    @{ Old = \"94\u00d747=\"; New = \"61\u00d785=\" },\n    @{ Old = \"29\u00d782=\"; New = \"70\u00d761=\" },\n    @{ Old = \"81\u00d743=\"; New = \"39\u00d790=\" },\n    @{ Old = \"72\u00d798=\"; New = \"48\u00d745=\" },\n    @{ Old = \"56\u00d760=\"; New = \"55\u00d773=\" },\n    @{ Old = \"18\u00d768=\"; New = \"53\u00d737=\" },\n    @{ Old = \"71\u00d751=\"; New = \"73\u00d795=\" },\n    @{ Old = \"52\u00d738=\"; New = \"21\u00d780=\" },\n    @{ Old = \"81\u00d780=\"; New = \"35\u00d725=\" },\n    @{ Old = \"84\u00d779=\"; New = \"31\u00d757=\" },\n    @{ Old = \"61\u00d743=\"; New = \"54\u00d783=\" },\n    @{ Old = \"34\u00d778=\"; New = \"73\u00d735=\" },\n    @{ Old = \"81\u00d737=\"; New = \"20\u00d731=\" },\n    @{ Old = \"84\u00d736=\"; New = \"78\u00d723=\" },\n    @{ Old = \"85\u00d764=\"; New = \"20\u00d712=\" }\n)\n\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, $wdReplaceAll)\n}\n"}
